# Generate Report for Handoff
#
# For the six "Ready for handoff" rows (7, 9, 10, 12, 13, 14) in the
# zh-cn / de-de / Overview sheets, refresh the handoff timestamps and
# mark the Priority column as a "handoff type" ("ht") record.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = 7,9,10,12,13,14

foreach ($r in $rows) {
    # Overview!G and de-de!H both tracked the same "Latest HO Xliff
    # Generate Date" timestamp for these source files - bump both to
    # the newly generated time.
    $wsOverview.Range("G$r").Value = "2016-09-04 08:24:00"
    $wsDeDe.Range("H$r").Value = "2016-09-04 08:24:00"

    # zh-cn!H tracked its own "Latest Handoff Datetime" for the same
    # handoff pass.
    $wsZhCn.Range("H$r").Value = "2016-09-04 08:23:54"

    # Priority is now flagged as "ht" (handoff type) for both locales.
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
